$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new weekly rows right before the current row 26, pushing the
# existing rows 26:35 down to 28:37 (dimension grows from T35 to T37).
$ws.Rows("26:27").Insert()

# New row 26 - Especial / $/caja 18 kilos granel / Región de O'Higgins
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44651
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100104
$ws.Range("H26").Value = "Frutos de pepita"
$ws.Range("I26").Value = 100104003
$ws.Range("J26").Value = "Membrillo"
$ws.Range("K26").Value = "Champion"
$ws.Range("L26").Value = "Especial"
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 12600
$ws.Range("O26").Value = 12600
$ws.Range("P26").Value = 12600
$ws.Range("Q26").Value = "`$/caja 18 kilos granel"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 700
$ws.Range("T26").Value = 18

# New row 27 - Primera / $/caja 18 kilos granel / Región de O'Higgins
$ws.Range("A27").Value = 9
$ws.Range("B27").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44651
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100104
$ws.Range("H27").Value = "Frutos de pepita"
$ws.Range("I27").Value = 100104003
$ws.Range("J27").Value = "Membrillo"
$ws.Range("K27").Value = "Champion"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 250
$ws.Range("N27").Value = 10800
$ws.Range("O27").Value = 10800
$ws.Range("P27").Value = 10800
$ws.Range("Q27").Value = "`$/caja 18 kilos granel"
$ws.Range("R27").Value = "Región de O'Higgins"
$ws.Range("S27").Value = 600
$ws.Range("T27").Value = 18
